$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update model labels (shared strings) with new memory addresses
$ws.Range("A2").Value = "<keras.engine.sequential.Sequential object at 0x000001DD5D0FD8B0>"
$ws.Range("A3").Value = "<keras.engine.sequential.Sequential object at 0x000001DD4FE50400>"
$ws.Range("A4").Value = "<keras.engine.sequential.Sequential object at 0x000001DD5E719880>"

# Update existing accuracy/loss values
$ws.Range("B2").Value = 0.4985714256763458
$ws.Range("C2").Value = 0.6944359540939331

$ws.Range("B3").Value = 0.4982142746448517
$ws.Range("C3").Value = 0.6935915946960449

$ws.Range("B4").Value = 0.4979464411735535
$ws.Range("C4").Value = 0.6951512694358826

# Add new row 5 for an additional model
$ws.Range("A5").Value = "<keras.engine.sequential.Sequential object at 0x000001DDAF2746D0>"
$ws.Range("B5").Value = 0.5
$ws.Range("C5").Value = 0.7518377900123596
